$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.076.25"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "3.577.44"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "578.11"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").Value = "186.34"
$ws.Range("E6").Value = "  -4.76%  "
$ws.Range("D7").Value = "3.573.15"
$ws.Range("E7").Value = "  -1.90%  "
$ws.Range("D8").Value = "0.618"
$ws.Range("E8").Value = "  -4.61%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "0.183"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("E11").Value = "  -4.62%  "
$ws.Range("D12").Value = "55.06"
$ws.Range("E12").Value = "  -5.71%  "
$ws.Range("D13").Value = "'0.0000304"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "9.53"
$ws.Range("E14").Value = "  -4.95%  "
$ws.Range("D15").Value = "4.155.47"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").Value = "19.64"
$ws.Range("E16").Value = "  -3.96%  "
$ws.Range("D17").Value = "3.579.33"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").Value = "70.011.01"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "12.57"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("E21").Value = "  -3.36%  "
$ws.Range("D22").Value = "493.06"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "19.48"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").Value = "4.97"
$ws.Range("E24").Value = "  -5.32%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "96.75"
$ws.Range("E25").Value = "  +5.68%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "4.38"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").Value = "11.55"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("E28").Value = "  -7.11%  "
$ws.Range("D29").Value = "9.32"
$ws.Range("E29").Value = "  -3.10%  "
$ws.Range("D30").Value = "7.73"
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").Value = "31.63"
$ws.Range("E31").Value = "  -3.85%  "
$ws.Range("D32").Value = "'12.10"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D33").Value = "65.81"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("E34").Value = "  -6.77%  "
$ws.Range("D35").Value = "'574.50"
$ws.Range("E35").Value = "  -7.11%  "
$ws.Range("D36").Value = "3.24"
$ws.Range("E36").Value = "  +14.41%  "
$ws.Range("D37").Value = "0.415"
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("D38").Value = "38.79"
$ws.Range("E38").Value = "  -3.81%  "
$ws.Range("D40").Value = "0.0₃0790"
$ws.Range("E40").Value = "  -5.35%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "3.18"
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "3.45"
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("E43").Value = "  -9.69%  "
$ws.Range("D44").Value = "3.06"
$ws.Range("E44").Value = "  -4.27%  "
$ws.Range("D45").Value = "3.55"
$ws.Range("E45").Value = "  +6.83%  "
$ws.Range("D46").Value = "3.189.38"
$ws.Range("E46").Value = "  -4.52%  "
$ws.Range("D47").Value = "'0.0440"
$ws.Range("E47").Value = "  -4.20%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "9.64"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("B49").Value = "OceanProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D49").Value = "1.56"
$ws.Range("E49").Value = "  +29.54%  "
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("E51").Value = "  +0.06%  "
